$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2 through 36 changes from 45674 to 45675
for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 3).Value = 45675
}
